# leer estilos e insertarlos
# Normalize the header row text to lowercase, preserving formatting/styles.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "nombre"
$ws.Range("C1").Value = "edad"
